# edit.ps1 - applies the "address trimming bug" documentation update to
# vvs2_report.docx, per the commit:
#   "Fixed bug where addresses came from the db with spaces in the end.
#    This made address deletion not work"
#
# Net visible-content changes:
#  1. Small wording/punctuation fixes in two existing paragraphs.
#  2. Two new paragraphs (+ a trailing blank spacer paragraph) describing the
#     address-trim bug and its fix, inserted right after the "telefone
#     impossiveis" bug note.
#  3. Removal of the now-superseded generic "App nao faz trim de strings"
#     bug paragraph (its content is superseded/replaced by the new, more
#     specific paragraphs added in step 2).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Small wording fixes
# ---------------------------------------------------------------------

# "remover um address ,criar uma sale" -> "remover um address, criar uma sale"
$d.Content.Find.Execute(
    "remover um address ,criar uma sale", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "remover um address, criar uma sale", 2) | Out-Null

# "alguma funcionalidades" -> "algumas funcionalidades"
$d.Content.Find.Execute(
    "encontradas alguma funcionalidades", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "encontradas algumas funcionalidades", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the two new paragraphs (+ blank spacer) describing the bug
#    right after the "telefone impossiveis (... demasiados digitos)"
#    paragraph.
# ---------------------------------------------------------------------

$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("demasiados")) {
        $targetIdx = $i
        break
    }
}

$anchor = $d.Paragraphs.Item($targetIdx)

# -- new paragraph 1: bug description (red highlight, like its neighbour) --
$anchor.Range.InsertParagraphAfter() | Out-Null
$p1 = $d.Paragraphs.Item($targetIdx + 1)
$p1.Range.Text = "Os addresses são retornados pela bd com espaços em branco, no fim o que pode causar erros quando são feitas comparações da string de addresses vindos da bd"
$p1.Format.FirstLineIndent = 18.0
$p1.Format.Alignment = 3
$p1.Range.HighlightColorIndex = 6

# -- new paragraph 2: the fix ("Sol: ...") --
$p1.Range.InsertParagraphAfter() | Out-Null
$p2 = $d.Paragraphs.Item($targetIdx + 2)
$p2.Range.Text = "`tSol: fazer trim das strings que são recebidas das bases de dados, na class AdressRowDataGateway."
$p2.Format.FirstLineIndent = 18.0
$p2.Format.Alignment = 3

# -- new paragraph 3: blank spacer paragraph --
$p2.Range.InsertParagraphAfter() | Out-Null
$p3 = $d.Paragraphs.Item($targetIdx + 3)
$p3.Format.FirstLineIndent = 35.4
$p3.Format.Alignment = 3

# ---------------------------------------------------------------------
# 3. Remove the now-superseded "App nao faz trim de strings" paragraph
# ---------------------------------------------------------------------

$removeIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Contains("App n") -and $d.Paragraphs.Item($i).Range.Text.Contains("faz trim de strings")) {
        $removeIdx = $i
        break
    }
}
if ($removeIdx -gt 0) {
    $d.Paragraphs.Item($removeIdx).Range.Delete() | Out-Null
}
